# Applies the "Work examples from real-life projects" wording updates
# across slides 1-3 of the presentation, editing only the targeted
# substrings while preserving existing run formatting.

function Replace-InTextRange {
    param($tr, [string]$old, [string]$new)
    $full = $tr.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        throw "Substring not found: [$old]"
    }
    $chars = $tr.Characters($idx + 1, $old.Length)
    $chars.Text = $new
}

$p = $ppt.ActivePresentation

# --- Slide 1: title + subtitle ---
$s1 = $p.Slides.Item(1)

$trTitle1 = $s1.Shapes.Item(1).TextFrame.TextRange
Replace-InTextRange $trTitle1 `
    "Specification By Examples for " `
    "Specification By Examples (SBE's) for "

$trSub1 = $s1.Shapes.Item(2).TextFrame.TextRange
Replace-InTextRange $trSub1 `
    " is to elaborates various scenarios involved in Person to Person Payments as Specification by Examples." `
    " is to elaborate a couple of multiple-step P2P scenarios presenting a flavour of Specification By Examples derived collaboratively between Business Analyst, Developer, Tester, and Architect."

# --- Slide 2: "Share a lunch bill" scenario ---
$s2 = $p.Slides.Item(2)

$trTitle2 = $s2.Shapes.Item(1).TextFrame.TextRange
Replace-InTextRange $trTitle2 `
    "Scenario: Share a lunch bill with friend " `
    "Scenario 1: Share a lunch bill with a friend "

$trBody2 = $s2.Shapes.Item(2).TextFrame.TextRange
Replace-InTextRange $trBody2 `
    "THEN`tI'm presented with my friends name and contact details to select and pay" `
    "THEN`tI'm presented with my friends name and contact details to select and pay from"
Replace-InTextRange $trBody2 `
    "GIVEN `tI'm presented with my friends name and contact details to select and pay" `
    "GIVEN `tI'm presented with my friends name and contact details to select and pay from"

# --- Slide 3: "Request money" scenario ---
$s3 = $p.Slides.Item(3)

$trTitle3 = $s3.Shapes.Item(1).TextFrame.TextRange
Replace-InTextRange $trTitle3 `
    "Scenario: Request money from a friend" `
    "Scenario 2: Request money from a friend"

$trBody3 = $s3.Shapes.Item(2).TextFrame.TextRange
Replace-InTextRange $trBody3 `
    "WHEN`tI select the friends contact to request mony from " `
    "WHEN`tI select the friends contact to request money from "
Replace-InTextRange $trBody3 `
    "GIVEN `tI'm presented to enter the amount to pay in 2 decimal places" `
    "GIVEN `tI'm presented to enter the amount to request in 2 decimal places"
